$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user rows appended at the bottom of the table (rows 39-52).
$newUsers = @(
    @{ Row = 39; Name = "BIANCA VICENTE";      Email = "BIANCA.VICENTE@light.org.ph";     Branch = "MALOLOS";     Id = 904; Remarks = "LO1" },
    @{ Row = 40; Name = "MISTY PAGTALUNAN";     Email = "MISTY.PAGTALUNAN@light.org.ph";    Branch = "MALOLOS";     Id = 905; Remarks = "LO2" },
    @{ Row = 41; Name = "ALVIN MORENO";         Email = "ALVIN.MORENO@light.org.ph";        Branch = "MALOLOS";     Id = 906; Remarks = "LO3" },
    @{ Row = 42; Name = "EMALYN LLAGAS";        Email = "EMALYN.LLAGAS@light.org.ph";       Branch = "MALOLOS";     Id = 907; Remarks = "LO4" },
    @{ Row = 43; Name = "JOVIELYN RAYMUNDO";    Email = "JOVIELYN.RAYMUNDO@light.org.ph";   Branch = "MALOLOS";     Id = 908; Remarks = "LO5" },
    @{ Row = 44; Name = "ELIZABETH BARASIGAN";  Email = "ELIZABETH.BARASIGAN@light.org.ph"; Branch = "MALOLOS";     Id = 87;  Remarks = "UNIT-OIC" },
    @{ Row = 45; Name = "ROXANNE ROQUE";        Email = "ROXANNE.ROQUE@light.org.ph";       Branch = "MALOLOS";     Id = 914; Remarks = "LO1" },
    @{ Row = 46; Name = "LEMUEL SANPEDRO ";     Email = "LEMUEL.SANPEDRO @light.org.ph";    Branch = "MALOLOS";     Id = 915; Remarks = "LO2" },
    @{ Row = 47; Name = "MICHELLE CAPERAL ";    Email = "MICHELLE.CAPERAL @light.org.ph";   Branch = "MALOLOS";     Id = 916; Remarks = "LO3" },
    @{ Row = 48; Name = "MIZHELLE BUHAT";       Email = "MIZHELLE.BUHAT@light.org.ph";      Branch = "MALOLOS";     Id = 917; Remarks = "LO4" },
    @{ Row = 49; Name = "ALVIN CABANTUGAN";     Email = "ALVIN.CABANTUGAN@light.org.ph";    Branch = "MALOLOS";     Id = 918; Remarks = "LO5" },
    @{ Row = 50; Name = "JERRY BALAGAT";        Email = "JERRY.BALAGAT@light.org.ph";       Branch = "MALOLOS";     Id = 142; Remarks = "UNIT-OIC" },
    @{ Row = 51; Name = "ALBERT BASCO";         Email = "albert.basco@light.org.ph";        Branch = "MALOLOS";     Id = 32;  Remarks = "MANAGER"; Hyperlink = $true },
    @{ Row = 52; Name = "SUZETTE MADAYAG";      Email = "suzette.madayag@light.org.ph";     Branch = "MAIN OFFICE"; Id = 1;   Remarks = "MANAGER"; Hyperlink = $true }
)

foreach ($u in $newUsers) {
    $r = $u.Row
    $ws.Range("A$r").Value = $u.Name
    $ws.Range("B$r").Value = $u.Email
    $ws.Range("C$r").Value = $u.Branch
    $ws.Range("D$r").Value = $u.Id
    $ws.Range("E$r").Value = $u.Remarks
}

# Last two new entries carry a mailto hyperlink on their email cell.
$ws.Hyperlinks.Add($ws.Range("B51"), "mailto:albert.basco@light.org.ph")
$ws.Hyperlinks.Add($ws.Range("B52"), "mailto:suzette.madayag@light.org.ph")

# Keep the plain (non-hyperlink) look for the new rows, matching the rest of
# the sheet which carries no explicit cell style.
$ws.Range("A39:E52").Style = "Normal"

# Move the active selection to D2, matching the saved workbook state.
$ws.Range("D2").Select() | Out-Null
